$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) cells whose new value reads as a plain number need an
# explicit Text format first, otherwise Excel auto-converts the assigned
# string into a numeric value (losing the "thousands-dot" display style
# used throughout this sheet, e.g. "1.669.22").
$textCells = @("D5", "D6", "D7", "D8", "D9", "D10", "D11", "D13", "D14", "D16", "D17", "D19", "D21", "D22", "D26", "D29", "D30", "D31", "D32", "D33", "D36", "D37", "D39", "D40", "D42", "D44", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.240.46'
$ws.Range("E2").Value = '  -6.09%  '
$ws.Range("D3").Value = '1.669.22'
$ws.Range("E4").Value = '  +0.49%  '
$ws.Range("D5").Value = '218.04'
$ws.Range("E5").Value = '  -3.69%  '
$ws.Range("D6").Value = '0.5061'
$ws.Range("E6").Value = '  -12.53%  '
$ws.Range("D7").Value = '1.006'
$ws.Range("E7").Value = '  +0.46%  '
$ws.Range("D8").Value = '0.2658'
$ws.Range("E8").Value = '  -2.81%  '
$ws.Range("D9").Value = '0.06337'
$ws.Range("E9").Value = '  -4.20%  '
$ws.Range("D10").Value = '21.57'
$ws.Range("E10").Value = '  -6.95%  '
$ws.Range("D11").Value = '0.07364'
$ws.Range("E11").Value = '  -2.33%  '
$ws.Range("D12").Value = '1.671.90'
$ws.Range("E12").Value = '  -4.06%  '
$ws.Range("D13").Value = '4.536'
$ws.Range("E13").Value = '  -3.64%  '
$ws.Range("D14").Value = '0.5797'
$ws.Range("E14").Value = '  -3.77%  '
$ws.Range("D15").Value = '1.894.95'
$ws.Range("E15").Value = '  -4.15%  '
$ws.Range("D16").Value = '0.000008532'
$ws.Range("E16").Value = '  -2.60%  '
$ws.Range("D17").Value = '64.77'
$ws.Range("E17").Value = '  -13.29%  '
$ws.Range("D18").Value = '26.164.68'
$ws.Range("E18").Value = '  -6.33%  '
$ws.Range("D19").Value = '4.932'
$ws.Range("E19").Value = '  -7.25%  '
$ws.Range("E20").Value = '  +0.54%  '
$ws.Range("D21").Value = '10.83'
$ws.Range("E21").Value = '  -4.11%  '
$ws.Range("D22").Value = '189.36'
$ws.Range("E22").Value = '  -7.84%  '
$ws.Range("E23").Value = '  -6.58%  '
$ws.Range("E24").Value = '  +0.46%  '
$ws.Range("D26").Value = '7.687'
$ws.Range("E26").Value = '  -4.57%  '
$ws.Range("E27").Value = '  -5.17%  '
$ws.Range("E28").Value = '  -2.95%  '
$ws.Range("D29").Value = '0.05779'
$ws.Range("E29").Value = '  -5.78%  '
$ws.Range("D30").Value = '1.277'
$ws.Range("E30").Value = '  -7.91%  '
$ws.Range("D31").Value = '1.323'
$ws.Range("E31").Value = '  -5.05%  '
$ws.Range("D32").Value = '3.527'
$ws.Range("E32").Value = '  -5.58%  '
$ws.Range("D33").Value = '3.512'
$ws.Range("E33").Value = '  -6.10%  '
$ws.Range("E34").Value = '  -2.34%  '
$ws.Range("E35").Value = '  -2.41%  '
$ws.Range("D36").Value = '0.5972'
$ws.Range("E36").Value = '  -6.85%  '
$ws.Range("D37").Value = '2.361'
$ws.Range("E37").Value = '  -2.31%  '
$ws.Range("E38").Value = '  -2.73%  '
$ws.Range("D39").Value = '0.01611'
$ws.Range("E39").Value = '  -3.33%  '
$ws.Range("D40").Value = '6.004'
$ws.Range("E40").Value = '  -2.30%  '
$ws.Range("D41").Value = '1.072.45'
$ws.Range("E41").Value = '  -4.70%  '
$ws.Range("D42").Value = '0.8609'
$ws.Range("E42").Value = '  -1.60%  '
$ws.Range("E43").Value = '  +0.56%  '
$ws.Range("D44").Value = '99.53'
$ws.Range("E44").Value = '  -0.41%  '
$ws.Range("D45").Value = '1.817.75'
$ws.Range("E45").Value = '  -3.77%  '
$ws.Range("D46").Value = '0.00000000112'
$ws.Range("E46").Value = '  +2.88%  '
$ws.Range("D47").Value = '55.71'
$ws.Range("E47").Value = '  -6.09%  '
$ws.Range("D48").Value = '1.004'
$ws.Range("E48").Value = '  +0.48%  '
$ws.Range("D49").Value = '8.087'
$ws.Range("E49").Value = '  -2.17%  '
$ws.Range("D50").Value = '0.4300'
$ws.Range("E50").Value = '  -2.53%  '
$ws.Range("D51").Value = '0.05182'
$ws.Range("E51").Value = '  -3.60%  '
